$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.23829066666667
$ws.Range("H2").Value = 63.714872
$ws.Range("I2").Value = 0.9042366413687101
$ws.Range("J2").Value = 0.90423664136871
$ws.Range("M2").Value = 24.27461233333333
$ws.Range("N2").Value = 72.823837
$ws.Range("O2").Value = 0.540681551915145
$ws.Range("P2").Value = 0.540681551915145
$ws.Range("Q2").Value = 515.551272555985
$ws.Range("R2").Value = 4639.961453003863
$ws.Range("S2").Value = 0.4889040705537726
$ws.Range("T2").Value = 0.4889040705537725
$ws.Range("G3").Value = 21.23829066666667
$ws.Range("H3").Value = 63.714872
$ws.Range("I3").Value = 0.9042366413687101
$ws.Range("J3").Value = 0.90423664136871
$ws.Range("M3").Value = 1.291159666666667
$ws.Range("O3").Value = 0.0287586966480594
$ws.Range("P3").Value = 0.0287586966480594
$ws.Range("Q3").Value = 27.42202429774311
$ws.Range("R3").Value = 246.798218679688
$ws.Range("S3").Value = 0.02600466726718281
$ws.Range("T3").Value = 0.02600466726718281
$ws.Range("G4").Value = 21.23829066666667
$ws.Range("H4").Value = 63.714872
$ws.Range("I4").Value = 0.9042366413687101
$ws.Range("J4").Value = 0.90423664136871
$ws.Range("M4").Value = 1.899338333333333
$ws.Range("N4").Value = 5.698015
$ws.Range("O4").Value = 0.04230498858547889
$ws.Range("P4").Value = 0.04230498858547889
$ws.Range("Q4").Value = 40.33869959767556
$ws.Range("R4").Value = 363.04829637908
$ws.Range("S4").Value = 0.03825372079167505
$ws.Range("T4").Value = 0.03825372079167504
$ws.Range("G5").Value = 21.23829066666667
$ws.Range("H5").Value = 63.714872
$ws.Range("I5").Value = 0.9042366413687101
$ws.Range("J5").Value = 0.90423664136871
$ws.Range("M5").Value = 16.14987566666667
$ws.Range("N5").Value = 48.449627
$ws.Range("O5").Value = 0.3597149037350217
$ws.Range("P5").Value = 0.3597149037350217
$ws.Range("Q5").Value = 342.9957536391938
$ws.Range("R5").Value = 3086.961782752744
$ws.Range("S5").Value = 0.3252673964036249
$ws.Range("T5").Value = 0.3252673964036249
$ws.Range("G6").Value = 21.23829066666667
$ws.Range("H6").Value = 63.714872
$ws.Range("I6").Value = 0.9042366413687101
$ws.Range("J6").Value = 0.90423664136871
$ws.Range("M6").Value = 0.3088903333333333
$ws.Range("N6").Value = 0.926671
$ws.Range("O6").Value = 0.006880081234867635
$ws.Range("P6").Value = 0.006880081234867635
$ws.Range("Q6").Value = 6.560302683456889
$ws.Range("R6").Value = 59.042724151112
$ws.Range("S6").Value = 0.006221221548160598
$ws.Range("T6").Value = 0.006221221548160598
$ws.Range("G7").Value = 21.23829066666667
$ws.Range("H7").Value = 63.714872
$ws.Range("I7").Value = 0.9042366413687101
$ws.Range("J7").Value = 0.90423664136871
$ws.Range("M7").Value = 0.9724443333333334
$ws.Range("N7").Value = 2.917333
$ws.Range("O7").Value = 0.02165977788142728
$ws.Range("P7").Value = 0.02165977788142728
$ws.Range("Q7").Value = 20.65305540848622
$ws.Range("R7").Value = 185.877498676376
$ws.Range("S7").Value = 0.01958556480429408
$ws.Range("T7").Value = 0.01958556480429408
$ws.Range("I8").Value = 0.05937834432696559
$ws.Range("J8").Value = 0.05937834432696559
$ws.Range("M8").Value = 24.27461233333333
$ws.Range("N8").Value = 72.823837
$ws.Range("O8").Value = 0.540681551915145
$ws.Range("P8").Value = 0.540681551915145
$ws.Range("Q8").Value = 33.85461236529567
$ws.Range("R8").Value = 304.6915112876611
$ws.Range("S8").Value = 0.0321047753608556
$ws.Range("T8").Value = 0.0321047753608556
$ws.Range("I9").Value = 0.05937834432696559
$ws.Range("J9").Value = 0.05937834432696559
$ws.Range("M9").Value = 1.291159666666667
$ws.Range("O9").Value = 0.0287586966480594
$ws.Range("P9").Value = 0.0287586966480594
$ws.Range("Q9").Value = 1.800717120276334
$ws.Range("S9").Value = 0.001707643791963222
$ws.Range("T9").Value = 0.001707643791963222
$ws.Range("I10").Value = 0.05937834432696559
$ws.Range("J10").Value = 0.05937834432696559
$ws.Range("M10").Value = 1.899338333333333
$ws.Range("N10").Value = 5.698015
$ws.Range("O10").Value = 0.04230498858547889
$ws.Range("P10").Value = 0.04230498858547889
$ws.Range("Q10").Value = 2.648914105921667
$ws.Range("R10").Value = 23.840226953295
$ws.Range("S10").Value = 0.002512000178976914
$ws.Range("T10").Value = 0.002512000178976914
$ws.Range("I11").Value = 0.05937834432696559
$ws.Range("J11").Value = 0.05937834432696559
$ws.Range("M11").Value = 16.14987566666667
$ws.Range("N11").Value = 48.449627
$ws.Range("O11").Value = 0.3597149037350217
$ws.Range("P11").Value = 0.3597149037350217
$ws.Range("Q11").Value = 22.52344024839234
$ws.Range("R11").Value = 202.710962235531
$ws.Range("S11").Value = 0.0213592754135194
$ws.Range("T11").Value = 0.0213592754135194
$ws.Range("I12").Value = 0.05937834432696559
$ws.Range("J12").Value = 0.05937834432696559
$ws.Range("M12").Value = 0.3088903333333333
$ws.Range("N12").Value = 0.926671
$ws.Range("O12").Value = 0.006880081234867635
$ws.Range("P12").Value = 0.006880081234867635
$ws.Range("Q12").Value = 0.4307942122736668
$ws.Range("R12").Value = 3.877147910463001
$ws.Range("S12").Value = 0.000408527832561465
$ws.Range("T12").Value = 0.0004085278325614651
$ws.Range("I13").Value = 0.05937834432696559
$ws.Range("J13").Value = 0.05937834432696559
$ws.Range("M13").Value = 0.9724443333333334
$ws.Range("N13").Value = 2.917333
$ws.Range("O13").Value = 0.02165977788142728
$ws.Range("P13").Value = 0.02165977788142728
$ws.Range("Q13").Value = 1.356220461927667
$ws.Range("R13").Value = 12.205984157349
$ws.Range("S13").Value = 0.001286121749088983
$ws.Range("T13").Value = 0.001286121749088983
$ws.Range("G14").Value = 0.6246503333333333
$ws.Range("H14").Value = 1.873951
$ws.Range("I14").Value = 0.02659497076804196
$ws.Range("J14").Value = 0.02659497076804196
$ws.Range("M14").Value = 24.27461233333333
$ws.Range("N14").Value = 72.823837
$ws.Range("O14").Value = 0.540681551915145
$ws.Range("P14").Value = 0.540681551915145
$ws.Range("Q14").Value = 15.16314468555411
$ws.Range("R14").Value = 136.468302169987
$ws.Range("S14").Value = 0.01437941006800284
$ws.Range("T14").Value = 0.01437941006800284
$ws.Range("G15").Value = 0.6246503333333333
$ws.Range("H15").Value = 1.873951
$ws.Range("I15").Value = 0.02659497076804196
$ws.Range("J15").Value = 0.02659497076804196
$ws.Range("M15").Value = 1.291159666666667
$ws.Range("O15").Value = 0.0287586966480594
$ws.Range("P15").Value = 0.0287586966480594
$ws.Range("Q15").Value = 0.8065233161698889
$ws.Range("R15").Value = 7.258709845529
$ws.Range("S15").Value = 0.0007648366966821261
$ws.Range("T15").Value = 0.0007648366966821261
$ws.Range("G16").Value = 0.6246503333333333
$ws.Range("H16").Value = 1.873951
$ws.Range("I16").Value = 0.02659497076804196
$ws.Range("J16").Value = 0.02659497076804196
$ws.Range("M16").Value = 1.899338333333333
$ws.Range("N16").Value = 5.698015
$ws.Range("O16").Value = 0.04230498858547889
$ws.Range("P16").Value = 0.04230498858547889
$ws.Range("Q16").Value = 1.186422323029444
$ws.Range("R16").Value = 10.677800907265
$ws.Range("S16").Value = 0.00112509993477316
$ws.Range("T16").Value = 0.00112509993477316
$ws.Range("G17").Value = 0.6246503333333333
$ws.Range("H17").Value = 1.873951
$ws.Range("I17").Value = 0.02659497076804196
$ws.Range("J17").Value = 0.02659497076804196
$ws.Range("M17").Value = 16.14987566666667
$ws.Range("N17").Value = 48.449627
$ws.Range("O17").Value = 0.3597149037350217
$ws.Range("P17").Value = 0.3597149037350217
$ws.Range("Q17").Value = 10.08802521847522
$ws.Range("R17").Value = 90.792226966277
$ws.Range("S17").Value = 0.009566607349661931
$ws.Range("T17").Value = 0.009566607349661931
$ws.Range("G18").Value = 0.6246503333333333
$ws.Range("H18").Value = 1.873951
$ws.Range("I18").Value = 0.02659497076804196
$ws.Range("J18").Value = 0.02659497076804196
$ws.Range("M18").Value = 0.3088903333333333
$ws.Range("N18").Value = 0.926671
$ws.Range("O18").Value = 0.006880081234867635
$ws.Range("P18").Value = 0.006880081234867635
$ws.Range("Q18").Value = 0.1929484496801111
$ws.Range("R18").Value = 1.736536047121
$ws.Range("S18").Value = 0.0001829755593230588
$ws.Range("T18").Value = 0.0001829755593230588
$ws.Range("G19").Value = 0.6246503333333333
$ws.Range("H19").Value = 1.873951
$ws.Range("I19").Value = 0.02659497076804196
$ws.Range("J19").Value = 0.02659497076804196
$ws.Range("M19").Value = 0.9724443333333334
$ws.Range("N19").Value = 2.917333
$ws.Range("O19").Value = 0.02165977788142728
$ws.Range("P19").Value = 0.02165977788142728
$ws.Range("Q19").Value = 0.6074376769647778
$ws.Range("R19").Value = 5.466939092683
$ws.Range("S19").Value = 0.0005760411595988404
$ws.Range("T19").Value = 0.0005760411595988404
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 0.6666666666666666
$ws.Range("G20").Value = 0.229944
$ws.Range("H20").Value = 0.689832
$ws.Range("I20").Value = 0.009790043536282392
$ws.Range("J20").Value = 0.009790043536282392
$ws.Range("M20").Value = 24.27461233333333
$ws.Range("N20").Value = 72.823837
$ws.Range("O20").Value = 0.540681551915145
$ws.Range("P20").Value = 0.540681551915145
$ws.Range("Q20").Value = 5.581801458376001
$ws.Range("R20").Value = 50.236213125384
$ws.Range("S20").Value = 0.005293295932513997
$ws.Range("T20").Value = 0.005293295932513997
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 0.6666666666666666
$ws.Range("G21").Value = 0.229944
$ws.Range("H21").Value = 0.689832
$ws.Range("I21").Value = 0.009790043536282392
$ws.Range("J21").Value = 0.009790043536282392
$ws.Range("M21").Value = 1.291159666666667
$ws.Range("O21").Value = 0.0287586966480594
$ws.Range("P21").Value = 0.0287586966480594
$ws.Range("Q21").Value = 0.296894418392
$ws.Range("R21").Value = 2.672049765528
$ws.Range("S21").Value = 0.00028154889223124
$ws.Range("T21").Value = 0.00028154889223124
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 0.6666666666666666
$ws.Range("G22").Value = 0.229944
$ws.Range("H22").Value = 0.689832
$ws.Range("I22").Value = 0.009790043536282392
$ws.Range("J22").Value = 0.009790043536282392
$ws.Range("M22").Value = 1.899338333333333
$ws.Range("N22").Value = 5.698015
$ws.Range("O22").Value = 0.04230498858547889
$ws.Range("P22").Value = 0.04230498858547889
$ws.Range("Q22").Value = 0.43674145372
$ws.Range("R22").Value = 3.93067308348
$ws.Range("S22").Value = 0.000414167680053768
$ws.Range("T22").Value = 0.000414167680053768
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 0.6666666666666666
$ws.Range("G23").Value = 0.229944
$ws.Range("H23").Value = 0.689832
$ws.Range("I23").Value = 0.009790043536282392
$ws.Range("J23").Value = 0.009790043536282392
$ws.Range("M23").Value = 16.14987566666667
$ws.Range("N23").Value = 48.449627
$ws.Range("O23").Value = 0.3597149037350217
$ws.Range("P23").Value = 0.3597149037350217
$ws.Range("Q23").Value = 3.713567010296
$ws.Range("R23").Value = 33.422103092664
$ws.Range("S23").Value = 0.003521624568215493
$ws.Range("T23").Value = 0.003521624568215493
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 0.6666666666666666
$ws.Range("G24").Value = 0.229944
$ws.Range("H24").Value = 0.689832
$ws.Range("I24").Value = 0.009790043536282392
$ws.Range("J24").Value = 0.009790043536282392
$ws.Range("M24").Value = 0.3088903333333333
$ws.Range("N24").Value = 0.926671
$ws.Range("O24").Value = 0.006880081234867635
$ws.Range("P24").Value = 0.006880081234867635
$ws.Range("Q24").Value = 0.071027478808
$ws.Range("R24").Value = 0.639247309272
$ws.Range("S24").Value = 0.00006735629482251366
$ws.Range("T24").Value = 0.00006735629482251368
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 0.6666666666666666
$ws.Range("G25").Value = 0.229944
$ws.Range("H25").Value = 0.689832
$ws.Range("I25").Value = 0.009790043536282392
$ws.Range("J25").Value = 0.009790043536282392
$ws.Range("M25").Value = 0.9724443333333334
$ws.Range("N25").Value = 2.917333
$ws.Range("O25").Value = 0.02165977788142728
$ws.Range("P25").Value = 0.02165977788142728
$ws.Range("Q25").Value = 0.223607739784
$ws.Range("R25").Value = 2.012469658056
$ws.Range("S25").Value = 0.0002120501684453795
$ws.Range("T25").Value = 0.0002120501684453795
